$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: drop "semester" and "jam_kompen"/"jam_kompen_selesai", add "periode_tahun".
# New order: nomor_induk, nama, periode_tahun, jam_alpha, prodi, (blank, styled)
$ws.Range("A1").Value = "nomor_induk"
$ws.Range("B1").Value = "nama"
$ws.Range("C1").Value = "periode_tahun"
$ws.Range("D1").Value = "jam_alpha"
$ws.Range("E1").Value = "prodi"
$ws.Range("F1").Value = ""
$ws.Range("G1").Clear()

# Column width tweaks (dashboard/filter layout change)
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 19.333333333333332

# Move the selection to D2 (was D8 with full-column sqref)
[void]$ws.Range("D2").Select()
